# "removed the crossed out items"
#
# The slide contains a diagram (Stack / Heap boxes) with a "crossed out"
# duplicate copy of part of the diagram lower on the slide (Rectangle 30,
# the two red "X" connectors, Group 32 "Point" struct, Group 36 "int"
# struct, their linking curved/straight connectors). That crossed-out
# copy is removed entirely, and the remaining shapes that used to sit
# below it are shifted up to fill the freed vertical space.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Move the surviving shapes up by the vacated gap (202.50141732283464 pt
#     == 2571768 EMU) -----------------------------------------------------

$s.Shapes.Item("Rectangle 23").Top = 1785926 / 12700
$s.Shapes.Item("Group 14").Top = 3037644 / 12700
$s.Shapes.Item("Group 15").Top = 2714620 / 12700
$s.Shapes.Item("Group 19").Top = 1928802 / 12700
$s.Shapes.Item("Group 24").Top = 2214554 / 12700

$curvedConnector29 = $s.Shapes.Item("Curved Connector 29")
$curvedConnector29.Top = 2678901 / 12700
# best-effort: the authored XML also drops this connector's <a:stCxn>
# (its start no longer snaps to shape 14); the COM surface exposes
# ConnectorFormat.BeginDisconnect() for this but it is a no-op in this
# runtime, so the start connection site is left attached.
$curvedConnector29.ConnectorFormat.BeginDisconnect()

# --- Delete the crossed-out shapes entirely -----------------------------

$s.Shapes.Item("Rectangle 30").Delete()
$s.Shapes.Item("Group 32").Delete()
$s.Shapes.Item("Group 36").Delete()
$s.Shapes.Item("Curved Connector 40").Delete()
$s.Shapes.Item("Straight Connector 44").Delete()
$s.Shapes.Item("Straight Connector 45").Delete()
$s.Shapes.Item("Straight Connector 48").Delete()
$s.Shapes.Item("Straight Connector 50").Delete()
